$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $dataSheet)
$newSheet.Name = "best strategies"

$arr = New-Object 'object[,]' 50,13
$arr[0,0] = 5
$arr[0,1] = 1
$arr[0,2] = 2
$arr[0,3] = 6
$arr[0,4] = 137.375
$arr[0,5] = 294.375
$arr[0,7] = 7
$arr[0,8] = 1
$arr[0,9] = 2
$arr[0,10] = 4
$arr[0,11] = 119
$arr[0,12] = 243.5
$arr[1,0] = 7
$arr[1,1] = 3
$arr[1,2] = 6
$arr[1,3] = 7
$arr[1,4] = 142.75
$arr[1,5] = 299.375
$arr[1,7] = 7
$arr[1,8] = 1
$arr[1,9] = 6
$arr[1,10] = 4
$arr[1,11] = 131.875
$arr[1,12] = 278.375
$arr[2,0] = 7
$arr[2,1] = 1
$arr[2,2] = 6
$arr[2,3] = 4
$arr[2,4] = 123.25
$arr[2,5] = 267.875
$arr[2,7] = 7
$arr[2,8] = 2
$arr[2,9] = 5
$arr[2,10] = 8
$arr[2,11] = 135.5
$arr[2,12] = 287.875
$arr[3,0] = 7
$arr[3,1] = 3
$arr[3,2] = 4
$arr[3,3] = 6
$arr[3,4] = 133.5
$arr[3,5] = 284.875
$arr[3,7] = 6
$arr[3,8] = 1
$arr[3,9] = 3
$arr[3,10] = 5
$arr[3,11] = 128.75
$arr[3,12] = 276.125
$arr[4,0] = 7
$arr[4,1] = 1
$arr[4,2] = 3
$arr[4,3] = 4
$arr[4,4] = 143.375
$arr[4,5] = 283.5
$arr[4,7] = 7
$arr[4,8] = 1
$arr[4,9] = 3
$arr[4,10] = 4
$arr[4,11] = 152.875
$arr[4,12] = 297.5
$arr[5,0] = 7
$arr[5,1] = 2
$arr[5,2] = 3
$arr[5,3] = 3
$arr[5,4] = 139
$arr[5,5] = 280.125
$arr[5,7] = 7
$arr[5,8] = 3
$arr[5,9] = 4
$arr[5,10] = 6
$arr[5,11] = 133.375
$arr[5,12] = 286.5
$arr[6,0] = 4
$arr[6,1] = 1
$arr[6,2] = 2
$arr[6,3] = 5
$arr[6,4] = 136.75
$arr[6,5] = 293.625
$arr[6,7] = 7
$arr[6,8] = 2
$arr[6,9] = 3
$arr[6,10] = 3
$arr[6,11] = 144.75
$arr[6,12] = 292.375
$arr[7,0] = 7
$arr[7,1] = 1
$arr[7,2] = 2
$arr[7,3] = 4
$arr[7,4] = 129
$arr[7,5] = 271.25
$arr[7,7] = 6
$arr[7,8] = 1
$arr[7,9] = 5
$arr[7,10] = 3
$arr[7,11] = 137.125
$arr[7,12] = 278
$arr[8,0] = 6
$arr[8,1] = 1
$arr[8,2] = 5
$arr[8,3] = 3
$arr[8,4] = 138.25
$arr[8,5] = 286
$arr[8,7] = 6
$arr[8,8] = 2
$arr[8,9] = 4
$arr[8,10] = 5
$arr[8,11] = 143.25
$arr[8,12] = 300.75
$arr[9,0] = 6
$arr[9,1] = 2
$arr[9,2] = 4
$arr[9,3] = 5
$arr[9,4] = 142.125
$arr[9,5] = 294.25
$arr[9,7] = 7
$arr[9,8] = 1
$arr[9,9] = 4
$arr[9,10] = 7
$arr[9,11] = 141.125
$arr[9,12] = 286.5
$arr[10,0] = 7
$arr[10,1] = 2
$arr[10,2] = 4
$arr[10,3] = 6
$arr[10,4] = 138.875
$arr[10,5] = 293.5
$arr[10,7] = 7
$arr[10,8] = 3
$arr[10,9] = 6
$arr[10,10] = 7
$arr[10,11] = 132.625
$arr[10,12] = 284.5
$arr[11,0] = 7
$arr[11,1] = 1
$arr[11,2] = 4
$arr[11,3] = 7
$arr[11,4] = 135.375
$arr[11,5] = 294.875
$arr[11,7] = 5
$arr[11,8] = 1
$arr[11,9] = 2
$arr[11,10] = 6
$arr[11,11] = 146.5
$arr[11,12] = 309.25
$arr[12,0] = 6
$arr[12,1] = 1
$arr[12,2] = 3
$arr[12,3] = 5
$arr[12,4] = 139.75
$arr[12,5] = 289.5
$arr[12,7] = 6
$arr[12,8] = 2
$arr[12,9] = 3
$arr[12,10] = 6
$arr[12,11] = 134.875
$arr[12,12] = 284.5
$arr[13,0] = 7
$arr[13,1] = 2
$arr[13,2] = 5
$arr[13,3] = 8
$arr[13,4] = 139
$arr[13,5] = 300.375
$arr[13,7] = 7
$arr[13,8] = 2
$arr[13,9] = 4
$arr[13,10] = 6
$arr[13,11] = 141.375
$arr[13,12] = 295.25
$arr[14,0] = 6
$arr[14,1] = 2
$arr[14,2] = 3
$arr[14,3] = 6
$arr[14,4] = 136.875
$arr[14,5] = 297.75
$arr[14,7] = 4
$arr[14,8] = 1
$arr[14,9] = 2
$arr[14,10] = 5
$arr[14,11] = 145.5
$arr[14,12] = 307.75
$arr[15,0] = 7
$arr[15,1] = 2
$arr[15,2] = 3
$arr[15,3] = 5
$arr[15,4] = 125.375
$arr[15,5] = 276.625
$arr[15,7] = 7
$arr[15,8] = 2
$arr[15,9] = 2
$arr[15,10] = 7
$arr[15,11] = 147.625
$arr[15,12] = 301.125
$arr[16,0] = 7
$arr[16,1] = 3
$arr[16,2] = 3
$arr[16,3] = 8
$arr[16,4] = 134.875
$arr[16,5] = 290
$arr[16,7] = 6
$arr[16,8] = 2
$arr[16,9] = 4
$arr[16,10] = 7
$arr[16,11] = 133.375
$arr[16,12] = 285.625
$arr[17,0] = 5
$arr[17,1] = 1
$arr[17,2] = 3
$arr[17,3] = 4
$arr[17,4] = 130.125
$arr[17,5] = 278
$arr[17,7] = 6
$arr[17,8] = 1
$arr[17,9] = 2
$arr[17,10] = 7
$arr[17,11] = 132.25
$arr[17,12] = 280.625
$arr[18,0] = 5
$arr[18,1] = 2
$arr[18,2] = 3
$arr[18,3] = 7
$arr[18,4] = 135.875
$arr[18,5] = 297
$arr[18,7] = 7
$arr[18,8] = 1
$arr[18,9] = 5
$arr[18,10] = 5
$arr[18,11] = 135.5
$arr[18,12] = 291.75
$arr[19,0] = 5
$arr[19,1] = 2
$arr[19,2] = 5
$arr[19,3] = 6
$arr[19,4] = 141.75
$arr[19,5] = 300.875
$arr[19,7] = 7
$arr[19,8] = 2
$arr[19,9] = 5
$arr[19,10] = 3
$arr[19,11] = 147.125
$arr[19,12] = 298.25
$arr[20,0] = 6
$arr[20,1] = 1
$arr[20,2] = 2
$arr[20,3] = 7
$arr[20,4] = 132.5
$arr[20,5] = 275.25
$arr[20,7] = 5
$arr[20,8] = 1
$arr[20,9] = 3
$arr[20,10] = 4
$arr[20,11] = 141.375
$arr[20,12] = 292.25
$arr[21,0] = 5
$arr[21,1] = 1
$arr[21,2] = 2
$arr[21,3] = 7
$arr[21,4] = 134.5
$arr[21,5] = 289
$arr[21,7] = 5
$arr[21,8] = 1
$arr[21,9] = 2
$arr[21,10] = 7
$arr[21,11] = 138.25
$arr[21,12] = 286.5
$arr[22,0] = 6
$arr[22,1] = 2
$arr[22,2] = 4
$arr[22,3] = 3
$arr[22,4] = 132.375
$arr[22,5] = 277.5
$arr[22,7] = 6
$arr[22,8] = 2
$arr[22,9] = 4
$arr[22,10] = 3
$arr[22,11] = 148.75
$arr[22,12] = 293.875
$arr[23,0] = 6
$arr[23,1] = 2
$arr[23,2] = 4
$arr[23,3] = 7
$arr[23,4] = 143.25
$arr[23,5] = 299.875
$arr[23,7] = 5
$arr[23,8] = 2
$arr[23,9] = 5
$arr[23,10] = 6
$arr[23,11] = 132
$arr[23,12] = 286
$arr[24,0] = 7
$arr[24,1] = 2
$arr[24,2] = 2
$arr[24,3] = 7
$arr[24,4] = 142.625
$arr[24,5] = 308.75
$arr[24,7] = 5
$arr[24,8] = 2
$arr[24,9] = 3
$arr[24,10] = 7
$arr[24,11] = 132.25
$arr[24,12] = 283.75
$arr[25,0] = 6
$arr[25,1] = 2
$arr[25,2] = 2
$arr[25,3] = 5
$arr[25,4] = 136.25
$arr[25,5] = 291.125
$arr[25,7] = 7
$arr[25,8] = 2
$arr[25,9] = 3
$arr[25,10] = 5
$arr[25,11] = 153.5
$arr[25,12] = 302.125
$arr[26,0] = 7
$arr[26,1] = 2
$arr[26,2] = 5
$arr[26,3] = 3
$arr[26,4] = 135.75
$arr[26,5] = 283.375
$arr[26,7] = 7
$arr[26,8] = 3
$arr[26,9] = 3
$arr[26,10] = 8
$arr[26,11] = 141.125
$arr[26,12] = 298
$arr[27,0] = 7
$arr[27,1] = 2
$arr[27,2] = 5
$arr[27,3] = 5
$arr[27,4] = 136.5
$arr[27,5] = 289.875
$arr[27,7] = 7
$arr[27,8] = 2
$arr[27,9] = 5
$arr[27,10] = 5
$arr[27,11] = 139
$arr[27,12] = 294.875
$arr[28,0] = 5
$arr[28,1] = 2
$arr[28,2] = 3
$arr[28,3] = 6
$arr[28,4] = 135.75
$arr[28,5] = 291.125
$arr[28,7] = 6
$arr[28,8] = 2
$arr[28,9] = 2
$arr[28,10] = 5
$arr[28,11] = 143.25
$arr[28,12] = 299.25
$arr[29,0] = 7
$arr[29,1] = 1
$arr[29,2] = 5
$arr[29,3] = 5
$arr[29,4] = 144
$arr[29,5] = 300.5
$arr[29,7] = 5
$arr[29,8] = 2
$arr[29,9] = 3
$arr[29,10] = 6
$arr[29,11] = 133.875
$arr[29,12] = 293.375
$arr[30,0] = 5
$arr[30,1] = 1
$arr[30,2] = 4
$arr[30,3] = 5
$arr[30,4] = 123
$arr[30,5] = 265.875
$arr[30,7] = 6
$arr[30,8] = 1
$arr[30,9] = 5
$arr[30,10] = 8
$arr[30,11] = 128.5
$arr[30,12] = 278.875
$arr[31,0] = 7
$arr[31,1] = 1
$arr[31,2] = 3
$arr[31,3] = 7
$arr[31,4] = 142.75
$arr[31,5] = 299.875
$arr[31,7] = 7
$arr[31,8] = 1
$arr[31,9] = 3
$arr[31,10] = 7
$arr[31,11] = 137.125
$arr[31,12] = 285.75
$arr[32,0] = 7
$arr[32,1] = 1
$arr[32,2] = 5
$arr[32,3] = 7
$arr[32,4] = 137.25
$arr[32,5] = 289.375
$arr[32,7] = 6
$arr[32,8] = 1
$arr[32,9] = 3
$arr[32,10] = 7
$arr[32,11] = 136.25
$arr[32,12] = 289.125
$arr[33,0] = 7
$arr[33,1] = 2
$arr[33,2] = 4
$arr[33,3] = 8
$arr[33,4] = 141.5
$arr[33,5] = 305.125
$arr[33,7] = 7
$arr[33,8] = 3
$arr[33,9] = 3
$arr[33,10] = 6
$arr[33,11] = 129.125
$arr[33,12] = 282.25
$arr[34,0] = 6
$arr[34,1] = 1
$arr[34,2] = 3
$arr[34,3] = 7
$arr[34,4] = 134
$arr[34,5] = 297.5
$arr[34,7] = 5
$arr[34,8] = 1
$arr[34,9] = 5
$arr[34,10] = 4
$arr[34,11] = 137.75
$arr[34,12] = 285.75
$arr[35,0] = 5
$arr[35,1] = 2
$arr[35,2] = 3
$arr[35,3] = 5
$arr[35,4] = 139.375
$arr[35,5] = 292.125
$arr[35,7] = 5
$arr[35,8] = 1
$arr[35,9] = 4
$arr[35,10] = 5
$arr[35,11] = 141.125
$arr[35,12] = 290.125
$arr[36,0] = 6
$arr[36,1] = 1
$arr[36,2] = 5
$arr[36,3] = 8
$arr[36,4] = 137.875
$arr[36,5] = 303.375
$arr[36,7] = 6
$arr[36,8] = 1
$arr[36,9] = 4
$arr[36,10] = 4
$arr[36,11] = 143.375
$arr[36,12] = 296.375
$arr[37,0] = 7
$arr[37,1] = 2
$arr[37,2] = 4
$arr[37,3] = 7
$arr[37,4] = 141.625
$arr[37,5] = 304.875
$arr[37,7] = 5
$arr[37,8] = 1
$arr[37,9] = 3
$arr[37,10] = 8
$arr[37,11] = 140.375
$arr[37,12] = 306.5
$arr[38,0] = 7
$arr[38,1] = 1
$arr[38,2] = 2
$arr[38,3] = 5
$arr[38,4] = 147
$arr[38,5] = 297.75
$arr[38,7] = 7
$arr[38,8] = 2
$arr[38,9] = 4
$arr[38,10] = 8
$arr[38,11] = 143.25
$arr[38,12] = 294.875
$arr[39,0] = 5
$arr[39,1] = 1
$arr[39,2] = 3
$arr[39,3] = 8
$arr[39,4] = 140
$arr[39,5] = 307
$arr[39,7] = 7
$arr[39,8] = 1
$arr[39,9] = 5
$arr[39,10] = 7
$arr[39,11] = 140
$arr[39,12] = 285.5
$arr[40,0] = 7
$arr[40,1] = 2
$arr[40,2] = 6
$arr[40,3] = 8
$arr[40,4] = 138.625
$arr[40,5] = 297.625
$arr[40,7] = 7
$arr[40,8] = 2
$arr[40,9] = 4
$arr[40,10] = 7
$arr[40,11] = 147.75
$arr[40,12] = 301.625
$arr[41,0] = 5
$arr[41,1] = 2
$arr[41,2] = 3
$arr[41,3] = 4
$arr[41,4] = 138.5
$arr[41,5] = 294.5
$arr[41,7] = 5
$arr[41,8] = 2
$arr[41,9] = 3
$arr[41,10] = 4
$arr[41,11] = 141.375
$arr[41,12] = 298.5
$arr[42,0] = 7
$arr[42,1] = 3
$arr[42,2] = 3
$arr[42,3] = 6
$arr[42,4] = 139.125
$arr[42,5] = 294.875
$arr[42,7] = 7
$arr[42,8] = 1
$arr[42,9] = 2
$arr[42,10] = 5
$arr[42,11] = 139.25
$arr[42,12] = 297.5
$arr[43,0] = 5
$arr[43,1] = 1
$arr[43,2] = 5
$arr[43,3] = 4
$arr[43,4] = 136.625
$arr[43,5] = 280
$arr[43,7] = 7
$arr[43,8] = 2
$arr[43,9] = 6
$arr[43,10] = 8
$arr[43,11] = 141
$arr[43,12] = 302.5
$arr[44,0] = 6
$arr[44,1] = 1
$arr[44,2] = 4
$arr[44,3] = 4
$arr[44,4] = 140.5
$arr[44,5] = 301.875
$arr[44,7] = 5
$arr[44,8] = 2
$arr[44,9] = 3
$arr[44,10] = 5
$arr[44,11] = 135.375
$arr[44,12] = 294.375
$arr[45,0] = 5
$arr[45,1] = 1
$arr[45,2] = 3
$arr[45,3] = 5
$arr[45,4] = 131.75
$arr[45,5] = 284.25
$arr[45,7] = 5
$arr[45,8] = 1
$arr[45,9] = 4
$arr[45,10] = 4
$arr[45,11] = 149.875
$arr[45,12] = 289.375
$arr[46,0] = 5
$arr[46,1] = 1
$arr[46,2] = 4
$arr[46,3] = 4
$arr[46,4] = 143.625
$arr[46,5] = 289.375
$arr[46,7] = 5
$arr[46,8] = 1
$arr[46,9] = 2
$arr[46,10] = 4
$arr[46,11] = 145.875
$arr[46,12] = 297.375
$arr[47,0] = 7
$arr[47,1] = 3
$arr[47,2] = 2
$arr[47,3] = 5
$arr[47,4] = 135.625
$arr[47,5] = 280
$arr[47,7] = 7
$arr[47,8] = 3
$arr[47,9] = 2
$arr[47,10] = 5
$arr[47,11] = 143.375
$arr[47,12] = 288.875
$arr[48,0] = 5
$arr[48,1] = 1
$arr[48,2] = 2
$arr[48,3] = 4
$arr[48,4] = 140.375
$arr[48,5] = 299
$arr[48,7] = 6
$arr[48,8] = 1
$arr[48,9] = 5
$arr[48,10] = 6
$arr[48,11] = 139.125
$arr[48,12] = 295.25
$arr[49,0] = 6
$arr[49,1] = 1
$arr[49,2] = 5
$arr[49,3] = 6
$arr[49,4] = 141.25
$arr[49,5] = 300.375
$arr[49,7] = 5
$arr[49,8] = 1
$arr[49,9] = 3
$arr[49,10] = 5
$arr[49,11] = 151.375
$arr[49,12] = 308.25

$newSheet.Range("A2:M51").Value = $arr

$newSheet.Range("E54").Formula = "=AVERAGE(E2:E53)"
$newSheet.Range("F54").Formula = "=AVERAGE(F2:F53)"
$newSheet.Range("L54").Formula = "=AVERAGE(L2:L53)"
$newSheet.Range("M54").Formula = "=AVERAGE(M2:M53)"

$newSheet.Range("L54").Select()

$wb.BreakLink("forceI.xlsx")

